$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrección de error behavior en los test case:
# actualiza usuario y números de siniestro usados en el caso de prueba

$ws.Range("C3").Value = "apellegrini"

# Las celdas de columna F guardan el número de siniestro como texto
# (con espacios finales). Se usa el prefijo de apóstrofe para forzar
# texto y conservar el formato/estilo original de la celda.
$ws.Range("F7").Value = "'1120170200969  "
$ws.Range("F6").Value = "'1220170301466  "
$ws.Range("F5").Value = "'0420172010228  "
$ws.Range("F3").Value = "'1120194100448  "
$ws.Range("F2").Value = "'1220194200694   "
$ws.Range("F4").Value = "'0420194406900"

# Actualiza la celda seleccionada activa
$ws.Range("E5").Select()
